$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D27").Value = "ACL 2023 Review"
$ws.Range("E27").Value = "https://tech.scatterlab.co.kr/acl2023-review/"

$ws.Range("D51").Value = "누적 방문수 500만 돌파, 그리고 IT 분야 크리에이터 선정"
$ws.Range("E51").Value = "https://bskyvision.com/entry/%EB%88%84%EC%A0%81-%EB%B0%A9%EB%AC%B8%EC%88%98-500%EB%A7%8C-%EB%8F%8C%ED%8C%8C-%EA%B7%B8%EB%A6%AC%EA%B3%A0-IT-%EB%B6%84%EC%95%BC-%ED%81%AC%EB%A6%AC%EC%97%90%EC%9D%B4%ED%84%B0-%EC%84%A0%EC%A0%95"
